$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.0982
$ws.Range("A6").Value = -21.61040000000001
$ws.Range("A7").Value = -21.57650000000001
$ws.Range("A16").Value = -20.28639999999999
$ws.Range("A20").Value = -22.23730000000002
